$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "ALC"; Row = 64; Sets = @{ "H" = 142862050; "I" = 500002000; "J" = 6068; "K" = 500002000; "L" = 6068; "M" = -500001752; "N" = -6564 }; Deletes = @() },
    @{ Sheet = "ALC"; Row = 67; Sets = @{ "H" = 142862050; "I" = 500002000; "J" = 6068; "K" = 500002000; "L" = 6068; "M" = -500001142; "N" = -7784 }; Deletes = @() },
    @{ Sheet = "ALC"; Row = 74; Sets = @{ "H" = 2551.5; "I" = 2551.5; "J" = 0; "K" = 2551.5; "L" = 0; "N" = -1615.5 }; Deletes = @("M") },
    @{ Sheet = "ALC"; Row = 77; Sets = @{ "H" = 2551.5; "I" = 2551.5; "J" = 0; "K" = 12757.5; "L" = 0; "N" = -8077.5 }; Deletes = @("M") },
    @{ Sheet = "ALC"; Row = 100; Sets = @{ "H" = 57701.11; "I" = 101264; "J" = 3247.5; "K" = 101264; "L" = 3247.5; "M" = -100723; "N" = -4329.5 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 33; Sets = @{ "H" = 1013.5; "I" = 1013.5; "J" = 0; "K" = 1013.5; "L" = 0; "M" = -684.5 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 36; Sets = @{ "H" = 1160; "I" = 1160; "J" = 0; "K" = 1160; "L" = 0; "M" = -814 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 61; Sets = @{ "H" = 4338; "I" = 4666.6665; "J" = 4009.3333; "K" = 4666.6665; "L" = 4009.3333; "M" = -4454.6665; "N" = -4433.3333 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 97; Sets = @{ "H" = 66713.125; "I" = 125456.125; "J" = 7970.125; "K" = 125456.125; "L" = 7970.125; "M" = -124960.125; "N" = -8962.125 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 102; Sets = @{ "H" = 1911.875; "I" = 0; "J" = 1911.875; "K" = 0; "M" = 1911.875; "N" = -5155.875 }; Deletes = @("L") },
    @{ Sheet = "ARM"; Row = 122; Sets = @{ "H" = 1680.875; "I" = 1525.4546; "J" = 2022.8; "K" = 4576.3638; "L" = 6068.4; "M" = -2126.3638; "N" = -10968.4 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 132; Sets = @{ "H" = 5388.14; "I" = 5982.478; "J" = 4986.0884; "K" = 17947.434; "L" = 14958.2652; "M" = -15417.434; "N" = -20018.2652 }; Deletes = @() },
    @{ Sheet = "ARM"; Row = 136; Sets = @{ "H" = 4338; "I" = 4666.6665; "J" = 4009.3333; "K" = 13999.9995; "L" = 12027.9999; "M" = -11449.9995; "N" = -17127.9999 }; Deletes = @() },
    @{ Sheet = "BSM"; Row = 86; Sets = @{ "H" = 2490.4348; "I" = 1987.5; "J" = 2758.6667; "K" = 1987.5; "L" = 2758.6667; "M" = -864.5; "N" = -5004.6667 }; Deletes = @() },
    @{ Sheet = "BSM"; Row = 89; Sets = @{ "H" = 2490.4348; "I" = 1987.5; "J" = 2758.6667; "K" = 9937.5; "L" = 13793.3335; "M" = -4321.5; "N" = -25025.3335 }; Deletes = @() },
    @{ Sheet = "BSM"; Row = 105; Sets = @{ "H" = 1964.4166; "I" = 1423.3334; "J" = 2505.5; "K" = 1423.3334; "L" = 2505.5; "M" = 323.6666; "N" = -5999.5 }; Deletes = @() },
    @{ Sheet = "CRP"; Row = 31; Sets = @{ "H" = 20858666; "I" = 55557324; "J" = 39471.832; "K" = 55557324; "L" = 39471.832; "M" = -55557029; "N" = -40061.832 }; Deletes = @() },
    @{ Sheet = "CRP"; Row = 34; Sets = @{ "H" = 20858666; "I" = 55557324; "J" = 39471.832; "K" = 55557324; "L" = 39471.832; "M" = -55557122; "N" = -39875.832 }; Deletes = @() },
    @{ Sheet = "CRP"; Row = 132; Sets = @{ "H" = 32264132; "I" = 71438580; "J" = 2820.5293; "K" = 214315740; "L" = 8461.5879; "M" = -214313210; "N" = -13521.5879 }; Deletes = @() },
    @{ Sheet = "GSM"; Row = 80; Sets = @{ "H" = 2911.1667; "I" = 1852.5; "J" = 3213.6428; "K" = 1852.5; "L" = 3213.6428; "M" = -854.5; "N" = -5209.6428 }; Deletes = @() },
    @{ Sheet = "GSM"; Row = 83; Sets = @{ "H" = 2911.1667; "I" = 1852.5; "J" = 3213.6428; "K" = 9262.5; "L" = 16068.214; "M" = -4270.5; "N" = -26052.214 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 68; Sets = @{ "H" = 1696.9697; "I" = 1577.6; "J" = 2070; "K" = 1577.6; "L" = 2070; "M" = -828.5999999999999; "N" = -3568 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 71; Sets = @{ "H" = 1696.9697; "I" = 1577.6; "J" = 2070; "K" = 7888; "L" = 10350; "M" = -4144; "N" = -17838 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 82; Sets = @{ "H" = 1226.2858; "I" = 980; "J" = 1363.1111; "K" = 980; "L" = 1363.1111; "M" = -619; "N" = -2085.1111 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 85; Sets = @{ "H" = 1226.2858; "I" = 980; "J" = 1363.1111; "K" = 980; "L" = 1363.1111; "M" = 268; "N" = -3859.1111 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 93; Sets = @{ "H" = 1326.1666; "I" = 775.6923; "J" = 1976.7273; "K" = 775.6923; "L" = 1976.7273; "M" = 472.3077; "N" = -4472.7273 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 100; Sets = @{ "H" = 1231.95; "I" = 1065.9333; "J" = 1730; "K" = 1065.9333; "L" = 1730; "M" = -524.9332999999999; "N" = -2812 }; Deletes = @() },
    @{ Sheet = "LTW"; Row = 132; Sets = @{ "H" = 13099.833; "I" = 19385.715; "J" = 4299.6; "K" = 58157.145; "L" = 12898.8; "M" = -55627.145; "N" = -17958.8 }; Deletes = @() },
    @{ Sheet = "WVR"; Row = 62; Sets = @{ "H" = 3001; "I" = 3001.3333; "J" = 3000; "K" = 3001.3333; "L" = 3000; "M" = -2377.3333; "N" = -4248 }; Deletes = @() },
    @{ Sheet = "WVR"; Row = 65; Sets = @{ "H" = 3001; "I" = 3001.3333; "J" = 3000; "K" = 15006.6665; "L" = 15000; "M" = -11886.6665; "N" = -21240 }; Deletes = @() },
    @{ Sheet = "WVR"; Row = 81; Sets = @{ "H" = 1610; "I" = 1610; "J" = 0; "K" = 3220; "L" = 0; "N" = -2159 }; Deletes = @("M") },
    @{ Sheet = "WVR"; Row = 84; Sets = @{ "H" = 1610; "I" = 1610; "J" = 0; "K" = 16100; "L" = 0; "N" = -10796 }; Deletes = @("M") },
    @{ Sheet = "WVR"; Row = 107; Sets = @{ "H" = 391.5; "I" = 388.66666; "J" = 400; "K" = 1165.99998; "L" = 1200; "M" = 754.0000199999999; "N" = -5040 }; Deletes = @() },
    @{ Sheet = "WVR"; Row = 136; Sets = @{ "H" = 4030.0833; "I" = 9702.909; "J" = 1534.04; "K" = 29108.727; "L" = 4602.12; "M" = -26558.727; "N" = -9702.119999999999 }; Deletes = @() },
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    foreach ($col in $e.Sets.Keys) {
        $addr = "$col$($e.Row)"
        $ws.Range($addr).Value = $e.Sets[$col]
    }
    foreach ($col in $e.Deletes) {
        $addr = "$col$($e.Row)"
        $ws.Range($addr).ClearContents()
    }
}

Write-Host "Applied $($edits.Count) row edits"